$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the numeric-looking Price/Volume columns stay stored as text (matching
# the source data's inline-string cells) instead of being auto-coerced to numbers.
$ws.Range("D2:E50").NumberFormat = "@"

$ws.Range("D2").Value = '255.92'
$ws.Range("E2").Value = '4.29%'
$ws.Range("D3").Value = '27.53'
$ws.Range("E3").Value = '-2.53%'
$ws.Range("D4").Value = '5.214'
$ws.Range("E4").Value = '-0.79%'
$ws.Range("D5").Value = '0.05923'
$ws.Range("E5").Value = '3.88%'
$ws.Range("D6").Value = '6.680'
$ws.Range("E6").Value = '0.71%'
$ws.Range("D7").Value = '0.8665'
$ws.Range("E7").Value = '1.84%'
$ws.Range("E8").Value = '14.72%'
$ws.Range("D9").Value = '0.1420'
$ws.Range("E9").Value = '3.43%'
$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").Value = '0.07183'
$ws.Range("E10").Value = '1.42%'
$ws.Range("B11").Value = 'BitrueCoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D11").Value = '0.03241'
$ws.Range("E11").Value = '2.43%'
$ws.Range("B12").Value = 'BitMartToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D12").Value = '0.09211'
$ws.Range("E12").Value = '-0.09%'
$ws.Range("B13").Value = 'BitForexToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D13").Value = '0.001546'
$ws.Range("E13").Value = '0.82%'
$ws.Range("B14").Value = 'One'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D14").Value = '0.0006104'
$ws.Range("E14").Value = '2.09%'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '0.005717'
$ws.Range("E15").Value = '-2.83%'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '3.484'
$ws.Range("E16").Value = '-0.20%'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '3.264'
$ws.Range("E17").Value = '1.90%'
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").Value = '2.206'
$ws.Range("E18").Value = '0.80%'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").Value = '0.3151'
$ws.Range("E19").Value = '-0.59%'
$ws.Range("B20").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C20").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D20").Value = '0.03632'
$ws.Range("E20").Value = '9.39%'
$ws.Range("D21").Value = '0.1309'
$ws.Range("E21").Value = '2.54%'
$ws.Range("D22").Value = '3.526'
$ws.Range("E22").Value = '-0.31%'
$ws.Range("D23").Value = '0.04177'
$ws.Range("E23").Value = '2.64%'
$ws.Range("E24").Value = '1.56%'
$ws.Range("E25").Value = '-0.14%'
$ws.Range("D26").Value = '0.004514'
$ws.Range("E26").Value = '8.78%'
$ws.Range("D27").Value = '0.0001201'
$ws.Range("E27").Value = '0.10%'
$ws.Range("D28").Value = '0.0001939'
$ws.Range("E28").Value = '33.86%'
$ws.Range("D40").Value = '0.03813'
$ws.Range("E40").Value = '0.40%'
$ws.Range("D41").Value = '0.005503'
$ws.Range("E41").Value = '47.24%'
$ws.Range("E42").Value = '3.29%'
$ws.Range("D43").Value = '0.001901'
$ws.Range("E43").Value = '-23.62%'
$ws.Range("D44").Value = '0.01067'
$ws.Range("E44").Value = '16.18%'
$ws.Range("D45").Value = '0.00005433'
$ws.Range("E45").Value = '2.69%'
$ws.Range("D46").Value = '0.00000000750'
$ws.Range("E46").Value = '0.07%'
$ws.Range("E47").Value = '3.98%'
$ws.Range("D48").Value = '0.002168'
$ws.Range("E48").Value = '-4.46%'
$ws.Range("D49").Value = '0.00002101'
$ws.Range("E49").Value = '0.07%'
$ws.Range("D50").Value = '0.0002001'
$ws.Range("E50").Value = '0.07%'
